$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the title heading ---
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Explore the land of the Pharaohs with Ancient Egypt Classics. Read our review and play this slot game for free today.</w:t></w:r>' + `
  '</w:p>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# --- 2. Remove the duplicate bold title paragraph near the end of the document ---
# Locate it precisely as a whole paragraph (the last occurrence of the title text)
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Play Ancient Egypt Classics for Free - Review") {
        $p.Range.Delete()
        break
    }
}

# --- 3. Replace the old meta-description-like sentence at the end with the new image prompt ---
$newImagePrompt = "Create a feature image for Ancient Egypt Classic that depicts a happy Maya warrior with glasses in a cartoon style. The warrior should be wearing traditional Mayan clothing and holding a golden scarab, representing the Wild symbol in the game. The background should show a temple with hieroglyphs and Egyptian gods, reflecting the theme of the game. The image should pop with bright colors and convey a fun and playful mood to attract potential players."

# Scope the Find to the last paragraph only, so the identical sentence that now also
# lives inside the new "Meta description" paragraph is left untouched.
$closingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$closingPara.Range.Find.Execute("Explore the land of the Pharaohs with Ancient Egypt Classics. Read our review and play this slot game for free today.", `
    $true, $false, $false, $false, $false, $true, 1, $false, $newImagePrompt, 2) | Out-Null
